# MAINT: Corrected VaR calculation and improved unittests
#
# 1) Add a 4th (empty) worksheet "Tabelle4" after the existing sheets.
# 2) On "Tabelle1": insert a "1 Day-Values" header above the existing
#    1-day VaR block (which moves down by one row and whose formulas
#    become row/column-absolute), then append "5 Day-Values" and
#    "10 Day-Values" blocks that reuse the same mean/std with an
#    n-day scaling factor.
# 3) On "Tabelle2": just move the selection.
# 4) Set the Tabelle1 page setup (paper size / orientation) and tidy
#    up the view (top-left cell / selection).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Tabelle1")
$ws2 = $wb.Worksheets.Item("Tabelle2")

# ---------------------------------------------------------------------
# 1) New empty sheet at the end of the workbook
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws4.Name = "Tabelle4"

# ---------------------------------------------------------------------
# 2) Tabelle1 - shift the old 1-day block down by two rows (11->13 .. 16->18),
#    leaving row 11 blank and row 12 holding the new section header, then
#    rewrite its formulas with absolute row/col refs where needed.
# ---------------------------------------------------------------------
$ws1.Range("A11:A12").EntireRow.Insert()

$ws1.Range("A12").Value = "1 Day-Values"
$ws1.Range("A12").Font.Bold = $true

$ws1.Range("B13").Formula = "=_xlfn.NORM.INV(1-D`$2, B`$9, B`$10)"
$ws1.Range("B14").Formula = "=C`$2*(B13+1)"
$ws1.Range("B15").Formula = "=C`$2-B14"
$ws1.Range("B16").Formula = "=_xlfn.NORM.INV(1-D`$3, B`$9, B`$10)"
$ws1.Range("B17").Formula = "=C`$2*(B16+1)"
$ws1.Range("B18").Formula = "=C`$2-B17"

# ---------------------------------------------------------------------
# 5 Day-Values block (rows 20-26)
# ---------------------------------------------------------------------
$ws1.Range("A20").Value = "5 Day-Values"
$ws1.Range("A20").Font.Bold = $true
$ws1.Range("B20").Value = 5
$ws1.Range("B20").Font.Bold = $true

$ws1.Range("A21").Value = "Min Return with 99% prob"
$ws1.Range("B21").Formula = "=_xlfn.NORM.INV(1-D`$2, B20*B`$9, SQRT(B20)*B`$10)"

$ws1.Range("A22").Value = "Value of Portfolio"
$ws1.Range("B22").Formula = "=C`$2*(B21+1)"

$ws1.Range("A23").Value = "value at Risk"
$ws1.Range("B23").Formula = "=C`$2-B22"

$ws1.Range("A24").Value = "Min Return with 95% prob"
$ws1.Range("B24").Formula = "=_xlfn.NORM.INV(1-D`$3, B20*B`$9, SQRT(B20)*B`$10)"

$ws1.Range("A25").Value = "Value of Portfolio"
$ws1.Range("B25").Formula = "=C`$2*(B24+1)"

$ws1.Range("A26").Value = "value at Risk"
$ws1.Range("B26").Formula = "=C`$2-B25"

# ---------------------------------------------------------------------
# 10 Day-Values block (rows 28-34)
# ---------------------------------------------------------------------
$ws1.Range("A28").Value = "10 Day-Values"
$ws1.Range("A28").Font.Bold = $true
$ws1.Range("B28").Value = 10
$ws1.Range("B28").Font.Bold = $true

# Documentation note lives in D20, but the shared-string table orders
# entries by first use -- "10 Day-Values" (A28) is referenced before this
# note in the original commit, so it is written here to match that order.
$ws1.Range("D20").Value = "See documentation in gg/powerline/doc/n_day_var.ipynb"

$ws1.Range("A29").Value = "Min Return with 99% prob"
$ws1.Range("B29").Formula = "=_xlfn.NORM.INV(1-D`$2, B28*B`$9, SQRT(B28)*B`$10)"

$ws1.Range("A30").Value = "Value of Portfolio"
$ws1.Range("B30").Formula = "=C`$2*(B29+1)"

$ws1.Range("A31").Value = "value at Risk"
$ws1.Range("B31").Formula = "=C`$2-B30"

$ws1.Range("A32").Value = "Min Return with 95% prob"
$ws1.Range("B32").Formula = "=_xlfn.NORM.INV(1-D`$3, B28*B`$9, SQRT(B28)*B`$10)"

$ws1.Range("A33").Value = "Value of Portfolio"
$ws1.Range("B33").Formula = "=C`$2*(B32+1)"

$ws1.Range("A34").Value = "value at Risk"
$ws1.Range("B34").Formula = "=C`$2-B33"

# ---------------------------------------------------------------------
# Page setup on Tabelle1
# ---------------------------------------------------------------------
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# Tabelle2 selection (done before re-activating Tabelle1 below, so the
# "selected tab" marker ends up back on Tabelle1, matching the source)
# ---------------------------------------------------------------------
$ws2.Range("A26").Select()

# ---------------------------------------------------------------------
# View tidy-up on Tabelle1 - re-activate it and scroll/select as needed
# ---------------------------------------------------------------------
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 7
$ws1.Range("C15").Select()
